$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (d=1)
$ws.Range("B2").Value = 0.05536738572306141
$ws.Range("C2").Value = -0.11034018169126
$ws.Range("D2").Value = 0.768606875789034
$ws.Range("E2").Value = -0.03477432211494006

# Row 3 (d=2)
$ws.Range("B3").Value = 0.1140651878936602
$ws.Range("C3").Value = -0.06539443422077537
$ws.Range("D3").Value = -0.01442959368523549
$ws.Range("E3").Value = 0.1049908797155545

# Row 4 (d=3)
$ws.Range("B4").Value = -0.07505793907031309
$ws.Range("C4").Value = 0.1062176505240186
$ws.Range("D4").Value = -0.1368056665725507
$ws.Range("E4").Value = 0.06839031063409196

# Row 5 (d=4)
$ws.Range("B5").Value = 0.1034705509471265
$ws.Range("C5").Value = -0.1463418316517127
$ws.Range("D5").Value = -0.06232977129326166
$ws.Range("E5").Value = 0.1265871775773947

# Row 6 (d=5)
$ws.Range("B6").Value = -1.414803524095609
$ws.Range("C6").Value = 11.50232248106203
$ws.Range("D6").Value = 12.70517222848191
$ws.Range("E6").Value = 25.86448165305455

# Row 7 (d=7)
$ws.Range("B7").Value = 98.73006621367118
$ws.Range("C7").Value = 98.33052178809814
$ws.Range("D7").Value = 98.24200989021298
$ws.Range("E7").Value = 97.5389890457102

# Row 8 (d=10)
$ws.Range("B8").Value = 95.64199836589123
$ws.Range("C8").Value = 94.46281468534568
$ws.Range("D8").Value = 95.55684902484158
$ws.Range("E8").Value = 95.80526163096499
